$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "63×54=3402" "70×28=1960"
Replace-Text "80×16=1280" "66×14=924"
Replace-Text "34×63=2142" "59×60=3540"
Replace-Text "77×19=1463" "16×14=224"
Replace-Text "26×86=2236" "35×88=3080"
Replace-Text "66×17=1122" "65×54=3510"
Replace-Text "70×88=6160" "61×23=1403"
Replace-Text "32×16=512" "22×50=1100"
Replace-Text "27×92=2484" "84×61=5124"
Replace-Text "46×34=1564" "44×60=2640"
Replace-Text "69×96=6624" "60×73=4380"
Replace-Text "90×62=5580" "26×34=884"
Replace-Text "74×90=6660" "96×12=1152"
Replace-Text "81×42=3402" "54×55=2970"
Replace-Text "63×89=5607" "54×45=2430"
Replace-Text "96×63=6048" "82×81=6642"
Replace-Text "28×16=448" "24×52=1248"
Replace-Text "99×54=5346" "20×25=500"
Replace-Text "64×61=3904" "41×89=3649"
Replace-Text "69×73=5037" "18×41=738"
Replace-Text "78×99=7722" "38×99=3762"
Replace-Text "62×38=2356" "37×94=3478"
Replace-Text "37×93=3441" "13×69=897"
Replace-Text "53×85=4505" "68×15=1020"
Replace-Text "93×56=5208" "45×29=1305"
